$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 965.05
$ws.Range("I41").Value = 1110.8
$ws.Range("J41").Value = 527.8
$ws.Range("K41").Value = 1110.8
$ws.Range("L41").Value = 527.8
$ws.Range("M41").Value = -670.8
$ws.Range("N41").Value = -1407.8
$ws.Range("H98").Value = 1099.8334
$ws.Range("I98").Value = 1061.7894
$ws.Range("J98").Value = 1244.4
$ws.Range("K98").Value = 1061.7894
$ws.Range("L98").Value = 1244.4
$ws.Range("M98").Value = 436.2106000000001
$ws.Range("N98").Value = -4240.4
$ws.Range("H113").Value = 2758
$ws.Range("I113").Value = 2642.8572
$ws.Range("J113").Value = 2873.1428
$ws.Range("K113").Value = 2642.8572
$ws.Range("L113").Value = 2873.1428
$ws.Range("M113").Value = 611.1428000000001
$ws.Range("N113").Value = -9381.1428
$ws.Range("H122").Value = 1099.8334
$ws.Range("I122").Value = 1061.7894
$ws.Range("J122").Value = 1244.4
$ws.Range("K122").Value = 3185.3682
$ws.Range("L122").Value = 3733.2
$ws.Range("M122").Value = -735.3681999999999
$ws.Range("N122").Value = -8633.200000000001
$ws.Range("H129").Value = 2428.0466
$ws.Range("I129").Value = 418.15384
$ws.Range("J129").Value = 3299
$ws.Range("K129").Value = 1254.46152
$ws.Range("L129").Value = 9897
$ws.Range("M129").Value = 3745.53848
$ws.Range("N129").Value = -19897
$ws.Range("H137").Value = 3925741.5
$ws.Range("I137").Value = 3715.24
$ws.Range("J137").Value = 7696920.5
$ws.Range("K137").Value = 11145.72
$ws.Range("L137").Value = 23090761.5
$ws.Range("M137").Value = -8595.719999999999
$ws.Range("N137").Value = -23095861.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1341.3636
$ws.Range("I2").Value = 1100
$ws.Range("J2").Value = 1431.875
$ws.Range("K2").Value = 1100
$ws.Range("L2").Value = 1431.875
$ws.Range("M2").Value = -987
$ws.Range("N2").Value = -1657.875
$ws.Range("H32").Value = 10763.695
$ws.Range("I32").Value = 11985.789
$ws.Range("J32").Value = 4958.75
$ws.Range("K32").Value = 11985.789
$ws.Range("L32").Value = 4958.75
$ws.Range("M32").Value = -11698.789
$ws.Range("N32").Value = -5532.75
$ws.Range("H61").Value = 13892291
$ws.Range("I61").Value = 22730114
$ws.Range("J61").Value = 4285.7144
$ws.Range("K61").Value = 22730114
$ws.Range("L61").Value = 4285.7144
$ws.Range("M61").Value = -22729902
$ws.Range("N61").Value = -4709.7144
$ws.Range("H116").Value = 1341.3636
$ws.Range("I116").Value = 1100
$ws.Range("J116").Value = 1431.875
$ws.Range("K116").Value = 1100
$ws.Range("L116").Value = 1431.875
$ws.Range("M116").Value = 1194
$ws.Range("N116").Value = -6019.875
$ws.Range("H122").Value = 5313.1924
$ws.Range("I122").Value = 5313.1924
$ws.Range("K122").Value = 15939.5772
$ws.Range("M122").Value = -13489.5772
$ws.Range("H136").Value = 13892291
$ws.Range("I136").Value = 22730114
$ws.Range("J136").Value = 4285.7144
$ws.Range("K136").Value = 68190342
$ws.Range("L136").Value = 12857.1432
$ws.Range("M136").Value = -68187792
$ws.Range("N136").Value = -17957.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1341.3636
$ws.Range("I3").Value = 1100
$ws.Range("J3").Value = 1431.875
$ws.Range("K3").Value = 1100
$ws.Range("L3").Value = 1431.875
$ws.Range("M3").Value = -986
$ws.Range("N3").Value = -1659.875
$ws.Range("H115").Value = 48990
$ws.Range("J115").Value = 48990
$ws.Range("L115").Value = 48990
$ws.Range("N115").Value = -52124
$ws.Range("H134").Value = 4092.7297
$ws.Range("I134").Value = 3825.761
$ws.Range("J134").Value = 4531.3213
$ws.Range("K134").Value = 11477.283
$ws.Range("L134").Value = 13593.9639
$ws.Range("M134").Value = -8942.282999999999
$ws.Range("N134").Value = -18663.9639

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 968.125
$ws.Range("I16").Value = 873.75
$ws.Range("J16").Value = 1062.5
$ws.Range("K16").Value = 873.75
$ws.Range("L16").Value = 1062.5
$ws.Range("M16").Value = -586.75
$ws.Range("N16").Value = -1636.5
$ws.Range("H31").Value = 15876963
$ws.Range("I31").Value = 3049.0557
$ws.Range("J31").Value = 111120450
$ws.Range("K31").Value = 3049.0557
$ws.Range("L31").Value = 111120450
$ws.Range("M31").Value = -2754.0557
$ws.Range("N31").Value = -111121040
$ws.Range("H34").Value = 15876963
$ws.Range("I34").Value = 3049.0557
$ws.Range("J34").Value = 111120450
$ws.Range("K34").Value = 3049.0557
$ws.Range("L34").Value = 111120450
$ws.Range("M34").Value = -2847.0557
$ws.Range("N34").Value = -111120854
$ws.Range("H94").Value = 3421.4285
$ws.Range("I94").Value = 1153.5714
$ws.Range("J94").Value = 7957.143
$ws.Range("K94").Value = 1153.5714
$ws.Range("L94").Value = 7957.143
$ws.Range("M94").Value = -702.5714
$ws.Range("N94").Value = -8859.143
$ws.Range("H113").Value = 968.125
$ws.Range("I113").Value = 873.75
$ws.Range("J113").Value = 1062.5
$ws.Range("K113").Value = 873.75
$ws.Range("L113").Value = 1062.5
$ws.Range("M113").Value = 1296.25
$ws.Range("N113").Value = -5402.5
$ws.Range("H132").Value = 2971.5925
$ws.Range("I132").Value = 2155.762
$ws.Range("K132").Value = 6467.286
$ws.Range("M132").Value = -3937.286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 18.869566
$ws.Range("I12").Value = 15.2
$ws.Range("J12").Value = 19.88889
$ws.Range("K12").Value = 45.59999999999999
$ws.Range("L12").Value = 59.66667
$ws.Range("M12").Value = 127.4
$ws.Range("N12").Value = -405.66667
$ws.Range("H17").Value = 800.1667
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 900.2
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 2700.6
$ws.Range("M17").Value = -731
$ws.Range("N17").Value = -3038.6
$ws.Range("H131").Value = 905.1772
$ws.Range("I131").Value = 470.81818
$ws.Range("J131").Value = 975.44116
$ws.Range("K131").Value = 1412.45454
$ws.Range("L131").Value = 2926.32348
$ws.Range("M131").Value = 3627.54546
$ws.Range("N131").Value = -13006.32348
$ws.Range("H132").Value = 1760.3077
$ws.Range("I132").Value = 1323
$ws.Range("J132").Value = 2460
$ws.Range("K132").Value = 11907
$ws.Range("L132").Value = 22140
$ws.Range("M132").Value = -9377
$ws.Range("N132").Value = -27200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17477.645
$ws.Range("I70").Value = 27022.54
$ws.Range("K70").Value = 27022.54
$ws.Range("M70").Value = -26752.54
$ws.Range("H73").Value = 17477.645
$ws.Range("I73").Value = 27022.54
$ws.Range("K73").Value = 27022.54
$ws.Range("M73").Value = -26086.54
$ws.Range("H80").Value = 24605842
$ws.Range("I80").Value = 30305758
$ws.Range("J80").Value = 8931071
$ws.Range("K80").Value = 30305758
$ws.Range("L80").Value = 8931071
$ws.Range("M80").Value = -30304760
$ws.Range("N80").Value = -8933067
$ws.Range("H83").Value = 24605842
$ws.Range("I83").Value = 30305758
$ws.Range("J83").Value = 8931071
$ws.Range("K83").Value = 151528790
$ws.Range("L83").Value = 44655355
$ws.Range("M83").Value = -151523798
$ws.Range("N83").Value = -44665339
$ws.Range("H102").Value = 2113.5386
$ws.Range("I102").Value = 2805.8696
$ws.Range("J102").Value = 1118.3125
$ws.Range("K102").Value = 2805.8696
$ws.Range("L102").Value = 1118.3125
$ws.Range("M102").Value = -1183.8696
$ws.Range("N102").Value = -4362.3125
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3705338
$ws.Range("I122").Value = 5556811
$ws.Range("K122").Value = 16670433
$ws.Range("M122").Value = -16667983
$ws.Range("H123").Value = 25321.666
$ws.Range("J123").Value = 25321.666
$ws.Range("L123").Value = 25321.666
$ws.Range("N123").Value = -30221.666
$ws.Range("H126").Value = 3668.0312
$ws.Range("I126").Value = 2155.6428
$ws.Range("K126").Value = 6466.928400000001
$ws.Range("M126").Value = -3996.928400000001
$ws.Range("H132").Value = 4317.8647
$ws.Range("I132").Value = 2773.1428
$ws.Range("J132").Value = 6345.3125
$ws.Range("K132").Value = 8319.428400000001
$ws.Range("L132").Value = 19035.9375
$ws.Range("M132").Value = -5789.428400000001
$ws.Range("N132").Value = -24095.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2137.5
$ws.Range("I16").Value = 2125
$ws.Range("J16").Value = 2275
$ws.Range("K16").Value = 2125
$ws.Range("L16").Value = 2275
$ws.Range("M16").Value = -1955
$ws.Range("N16").Value = -2615
$ws.Range("H46").Value = 1160
$ws.Range("I46").Value = 840
$ws.Range("J46").Value = 1337.7778
$ws.Range("K46").Value = 840
$ws.Range("L46").Value = 1337.7778
$ws.Range("M46").Value = -652
$ws.Range("N46").Value = -1713.7778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2225
$ws.Range("I81").Value = 1350
$ws.Range("J81").Value = 3100
$ws.Range("K81").Value = 2700
$ws.Range("L81").Value = 6200
$ws.Range("M81").Value = -1639
$ws.Range("N81").Value = -8322
$ws.Range("H84").Value = 2225
$ws.Range("I84").Value = 1350
$ws.Range("J84").Value = 3100
$ws.Range("K84").Value = 13500
$ws.Range("L84").Value = 31000
$ws.Range("M84").Value = -8196
$ws.Range("N84").Value = -41608
